$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")
$ws.Activate()

# --- Row 21 currently blank: seed its formatting from row 19 (B/C/D/E),
#     which -- at this point -- still carries the style indices (13/17/7/7)
#     that row 21 needs to end up with. Do this BEFORE row 19's own text
#     and style change below.
$ws.Range("B19:E19").Copy()
$ws.Range("B21:E21").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 21 new data ---
$ws.Range("B21").Value = 10
$ws.Range("C21").Value = "Site Estático Cadastro e Login - Local"
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 7

# --- Row 19: restyle C19 to match C18's look (font/border/alignment),
#     then replace its text with the new requirement.
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C19").Value = "Script do Banco de Dados"
$ws.Range("D19").Value = 13

# --- Row 20: swap in the text that used to sit in C19, keep its own style ---
$ws.Range("C20").Value = "Site Estático Dashboard (Gráfico com ChartJS) - Local"
$ws.Range("D20").Value = 21
$ws.Range("E20").Value = 10

# --- Simple size (story point) tweaks ---
$ws.Range("D13").Value = 8
$ws.Range("D14").Value = 8
$ws.Range("D18").Value = 8

# --- Selection moves to G18 ---
$ws.Range("G18").Select() | Out-Null
